# Improve and document logic of function parseRoundsTokensForPlayers
# Add a new "High card" example block (rows 17-19) to Sheet1, mirroring
# the existing FourOfAKind / FullHouse / Two pairs / Pair blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New header row (17): exponents used by the two data rows below ---
$ws.Range("B17").Value = 5
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0

# --- Row 18: "High card" labelled example ---
$ws.Range("A18").Value = "High card"
$ws.Range("B18").Value = 14
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 4
$ws.Range("F18").Value = 5
$ws.Range("G18").Formula = '=B18*POWER($J$1, $B$17) + C18*POWER($J$1, $C$17) + D18 *POWER($J$1,$D$17)+ E18 *POWER($J$1,$E$17)+F18 *POWER($J$1,$F$17)'

# --- Row 19: second example row for the High card block ---
$ws.Range("B19").Value = 13
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 8
$ws.Range("G19").Formula = '=B19*POWER($J$1, $B$17) + C19*POWER($J$1, $C$17) + D19 *POWER($J$1,$D$17)+ E19 *POWER($J$1,$E$17)+F19 *POWER($J$1,$F$17)'

# --- Formatting: reuse the thin-border "data row" style already used by
#     the other blocks (B2:F2 etc.) for the two new data rows. ---
$ws.Range("B2:F2").Copy()
$ws.Range("B18:F19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Match final selection shown in the workbook after the edit.
$ws.Range("G19").Select()
